# Update the "Förändrad" (Changed) date column (C) for rows 2-18
# from serial date 45206 (2023-10-07) to 45208 (2023-10-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45206) {
        $cell.Value = 45208
    }
}
